$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.922.60'
$ws.Range("D3").Value = '2.348.28'
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.667'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.44%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.98'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.599'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.33%  '
$ws.Range("E10").Value = '  +0.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '59.67'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '33.31'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.78%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.109'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.70%  '
$ws.Range("E14").Value = '  -1.41%  '
$ws.Range("D15").Value = '2.695.68'
$ws.Range("E15").Value = '  -0.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '16.15'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.906'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.70%  '
$ws.Range("D18").Value = '2.356.55'
$ws.Range("E18").Value = '  +0.21%  '
$ws.Range("D19").Value = '43.786.87'
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("E20").Value = '  +1.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '78.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.65'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '252.93'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.88%  '
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("E25").Value = '  +1.63%  '
$ws.Range("E26").Value = '  +2.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.50'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.44'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.31'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '175.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.11%  '
$ws.Range("E31").Value = '  -2.61%  '
$ws.Range("E32").Value = '  +0.95%  '
$ws.Range("E33").Value = '  -1.53%  '
$ws.Range("E34").Value = '  -0.93%  '
$ws.Range("E35").Value = '  -4.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.36'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.37%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.83'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.09%  '
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.39'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.03%  '
$ws.Range("B39").Value = 'THORChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.39'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.21%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.63'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +19.32%  '
$ws.Range("E41").Value = '  -3.38%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '64.75'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +16.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.17'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.96%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '18.91'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.29%  '
$ws.Range("E45").Value = '  -8.20%  '
$ws.Range("E46").Value = '  -1.97%  '
$ws.Range("E47").Value = '  -0.10%  '
$ws.Range("E48").Value = '  -1.92%  '
$ws.Range("E49").Value = '  -3.39%  '
$ws.Range("E50").Value = '  -1.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '98.48'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.21%  '
